$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.185.86'

$ws.Range("D3").Value = '3.324.05'
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'551.76"
$ws.Range("E5").Value = '  -0.04%  '

$ws.Range("D6").Value = "'172.54"
$ws.Range("E6").Value = '  -0.78%  '

$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = '  +1.90%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = '3.313.98'
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("E10").Value = '  +6.93%  '

$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("D12").Value = "'53.02"
$ws.Range("E12").Value = '  -1.43%  '

$ws.Range("D13").Value = "'0.0000277"
$ws.Range("E13").Value = '  +1.82%  '

$ws.Range("D14").Value = "'9.04"
$ws.Range("E14").Value = '  +0.16%  '

$ws.Range("D15").Value = '3.850.72'
$ws.Range("E15").Value = '  -0.64%  '

$ws.Range("E16").Value = '  +2.66%  '

$ws.Range("D17").Value = "'18.04"
$ws.Range("E17").Value = '  -1.35%  '

$ws.Range("D18").Value = '3.326.52'
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").Value = '64.058.17'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").Value = "'11.66"
$ws.Range("E20").Value = '  -0.83%  '

$ws.Range("D21").Value = "'0.980"
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").Value = "'450.53"
$ws.Range("E22").Value = '  +6.12%  '

$ws.Range("D23").Value = "'4.99"
$ws.Range("E23").Value = '  +3.60%  '

$ws.Range("E24").Value = '  -1.01%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = "'13.92"
$ws.Range("E25").Value = '  +6.70%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = "'86.97"
$ws.Range("E26").Value = '  +3.68%  '

$ws.Range("E27").Value = '  +1.38%  '

$ws.Range("D28").Value = "'10.54"
$ws.Range("E28").Value = '  -1.63%  '

$ws.Range("D29").Value = "'8.55"
$ws.Range("E29").Value = '  -0.54%  '

$ws.Range("D30").Value = "'30.78"
$ws.Range("E30").Value = '  +3.93%  '

$ws.Range("E31").Value = '  -2.27%  '

$ws.Range("D32").Value = "'62.68"
$ws.Range("E32").Value = '  +7.73%  '

$ws.Range("D33").Value = "'11.35"
$ws.Range("E33").Value = '  -0.39%  '

$ws.Range("D34").Value = "'569.11"
$ws.Range("E34").Value = '  -0.75%  '

$ws.Range("E35").Value = '  -0.87%  '

$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").Value = "'0.141"
$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").Value = "'3.51"
$ws.Range("E38").Value = '  +0.68%  '

$ws.Range("D39").Value = "'35.11"
$ws.Range("E39").Value = '  -1.09%  '

$ws.Range("D40").Value = "'0.365"
$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("E41").Value = '  -3.68%  '

$ws.Range("D42").Value = '3.057.08'
$ws.Range("E42").Value = '  -1.14%  '

$ws.Range("D43").Value = "'0.0411"

$ws.Range("E44").Value = '  -3.59%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = "'0.133"
$ws.Range("E46").Value = '  +2.83%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = "'3.14"
$ws.Range("E47").Value = '  -2.30%  '

$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("D49").Value = "'140.36"
$ws.Range("E49").Value = '  +4.45%  '

$ws.Range("E50").Value = '  -2.20%  '

$ws.Range("D51").Value = "'8.14"
$ws.Range("E51").Value = '  -0.42%  '
